$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply a temporary text number format to the ranges that receive new
# values so that numeric-looking strings (prices, percentages, hours)
# are stored as exact text rather than being auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

$ws.Range("D2").Value = "293.36"
$ws.Range("E2").Value = "2.38%"
$ws.Range("G2").Value = "23"
$ws.Range("D3").Value = "29.62"
$ws.Range("E3").Value = "3.74%"
$ws.Range("G3").Value = "23"
$ws.Range("D4").Value = "5.154"
$ws.Range("E4").Value = "1.89%"
$ws.Range("G4").Value = "23"
$ws.Range("D5").Value = "0.07137"
$ws.Range("E5").Value = "7.42%"
$ws.Range("G5").Value = "23"
$ws.Range("D6").Value = "7.512"
$ws.Range("E6").Value = "1.64%"
$ws.Range("G6").Value = "23"
$ws.Range("D7").Value = "3.623"
$ws.Range("E7").Value = "6.42%"
$ws.Range("G7").Value = "23"
$ws.Range("D8").Value = "1.403"
$ws.Range("E8").Value = "3.03%"
$ws.Range("G8").Value = "23"
$ws.Range("D9").Value = "0.9100"
$ws.Range("E9").Value = "-2.97%"
$ws.Range("G9").Value = "23"
$ws.Range("E10").Value = "3.35%"
$ws.Range("G10").Value = "23"
$ws.Range("D11").Value = "0.07647"
$ws.Range("E11").Value = "14.29%"
$ws.Range("G11").Value = "23"
$ws.Range("D12").Value = "0.07756"
$ws.Range("E12").Value = "2.65%"
$ws.Range("G12").Value = "23"
$ws.Range("D13").Value = "0.02936"
$ws.Range("E13").Value = "0.05%"
$ws.Range("G13").Value = "23"
$ws.Range("D14").Value = "0.09003"
$ws.Range("E14").Value = "0.20%"
$ws.Range("G14").Value = "23"
$ws.Range("D15").Value = "0.001602"
$ws.Range("E15").Value = "0.86%"
$ws.Range("G15").Value = "23"
$ws.Range("D16").Value = "0.0006556"
$ws.Range("E16").Value = "1.50%"
$ws.Range("G16").Value = "23"
$ws.Range("D17").Value = "0.006108"
$ws.Range("E17").Value = "-2.56%"
$ws.Range("G17").Value = "23"
$ws.Range("D18").Value = "3.494"
$ws.Range("E18").Value = "1.46%"
$ws.Range("G18").Value = "23"
$ws.Range("E19").Value = "-0.76%"
$ws.Range("G19").Value = "23"
$ws.Range("D20").Value = "0.3272"
$ws.Range("G20").Value = "23"
$ws.Range("D21").Value = "0.1367"
$ws.Range("E21").Value = "5.39%"
$ws.Range("G21").Value = "23"
$ws.Range("D22").Value = "4.025"
$ws.Range("E22").Value = "-1.56%"
$ws.Range("G22").Value = "23"
$ws.Range("G23").Value = "23"
$ws.Range("D24").Value = "0.04525"
$ws.Range("E24").Value = "0.60%"
$ws.Range("G24").Value = "23"
$ws.Range("D25").Value = "0.001210"
$ws.Range("E25").Value = "2.31%"
$ws.Range("G25").Value = "23"
$ws.Range("D26").Value = "0.004260"
$ws.Range("E26").Value = "2.78%"
$ws.Range("G26").Value = "23"
$ws.Range("D27").Value = "0.0001168"
$ws.Range("G27").Value = "23"
$ws.Range("D28").Value = "0.0001652"
$ws.Range("E28").Value = "2.20%"
$ws.Range("G28").Value = "23"
$ws.Range("G29").Value = "23"
$ws.Range("G30").Value = "23"
$ws.Range("G31").Value = "23"
$ws.Range("G32").Value = "23"
$ws.Range("G33").Value = "23"
$ws.Range("G34").Value = "23"
$ws.Range("G35").Value = "23"
$ws.Range("G36").Value = "23"
$ws.Range("G37").Value = "23"
$ws.Range("G38").Value = "23"
$ws.Range("G39").Value = "23"
$ws.Range("D40").Value = "0.04397"
$ws.Range("E40").Value = "4.70%"
$ws.Range("G40").Value = "23"
$ws.Range("D41").Value = "0.007018"
$ws.Range("E41").Value = "4.11%"
$ws.Range("G41").Value = "23"
$ws.Range("D42").Value = "0.1278"
$ws.Range("E42").Value = "2.32%"
$ws.Range("G42").Value = "23"
$ws.Range("D43").Value = "0.002207"
$ws.Range("G43").Value = "23"
$ws.Range("D44").Value = "0.01351"
$ws.Range("E44").Value = "12.38%"
$ws.Range("G44").Value = "23"
$ws.Range("D45").Value = "0.00005856"
$ws.Range("E45").Value = "4.16%"
$ws.Range("G45").Value = "23"
$ws.Range("G46").Value = "23"
$ws.Range("D47").Value = "0.01298"
$ws.Range("E47").Value = "-0.61%"
$ws.Range("G47").Value = "23"
$ws.Range("G48").Value = "23"
$ws.Range("G49").Value = "23"
$ws.Range("G50").Value = "23"
$ws.Range("G51").Value = "23"

# Restore the default "Normal" style so the cells keep the same
# (unstyled) appearance they had before this edit.
$ws.Range("D2:D51").Style = "Normal"
$ws.Range("E2:E51").Style = "Normal"
$ws.Range("G2:G51").Style = "Normal"
